$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scenario previously in row 5 (id 699) is removed; the scenario
# previously in row 6 (id 853) shifts up to become row 5. This is
# equivalent to deleting the entire row 5.
$ws.Rows.Item(5).Delete()
